# Update on 21st for single script
# - Column A "Run status" flips from Y to N for all data rows except the
#   first data row (row 2) on "MBA Standard Reports" and "Class Status".
# - The active sheet / selection moves from "Class Status" (last sheet)
#   to "MBA Standard Reports" (first sheet).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MBA Standard Reports")
$ws6 = $wb.Worksheets.Item("Class Status")

# Bulk-update column A (rows 3 through the last data row) from "Y" to "N".
$ws1.Range("A3:A121").Value = "N"
$ws6.Range("A3:A62").Value = "N"

# Restore the selection/scroll state on "Class Status" before moving away
# from it, then make "MBA Standard Reports" the active sheet with its new
# selection - matching the final view state in the workbook.
$ws6.Activate()
$ws6.Range("B76").Select()

$ws1.Activate()
$ws1.Range("C8").Select()
